# cierre de 10 de Nov 2021
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REMISIONES OCTUBRE  2021     ")

# D5 was "CANCELADA" -> now a blank/space placeholder (new shared string " ")
$ws.Range("D5").Value = " "

# New credit entries for Nov 2021 (rows 40-44), previously empty
$ws.Range("A40").Value = 44501
$ws.Range("D40").Value = "COMERCIO   CENTRAL "
$ws.Range("E40").Value = 259

$ws.Range("A41").Value = 44502
$ws.Range("D41").Value = "COMERCIO   CENTRAL "
$ws.Range("E41").Value = 8605

$ws.Range("A42").Value = 44503
$ws.Range("D42").Value = "COMERCIO   CENTRAL "
$ws.Range("E42").Value = 235

$ws.Range("A43").Value = 44505
$ws.Range("D43").Value = "OBRADOR"
$ws.Range("E43").Value = 1618

$ws.Range("A44").Value = 44505
$ws.Range("D44").Value = "COMERCIO   CENTRAL "
$ws.Range("E44").Value = 784

# Update the visible selection/cursor to reflect the new working cell
$ws.Activate()
$ws.Range("D4:D5").Select()
